$d = $word.ActiveDocument

$d.Content.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2) | Out-Null
$d.Content.Find.Execute("38÷6=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=7, 4", 2) | Out-Null
$d.Content.Find.Execute("39÷6=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷4=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=3, 2", 2) | Out-Null
$d.Content.Find.Execute("76÷9=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2) | Out-Null
$d.Content.Find.Execute("72÷6=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "15÷7=2, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷4=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "87÷2=43, 1", 2) | Out-Null
$d.Content.Find.Execute("20÷9=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=11, 4", 2) | Out-Null
$d.Content.Find.Execute("18÷9=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷9=8, 8", $true, $false, $false, $false, $false, $true, 1, $false, "34÷4=8, 2", 2) | Out-Null
$d.Content.Find.Execute("13÷5=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=2, 2", 2) | Out-Null
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "57÷7=8, 1", 2) | Out-Null
$d.Content.Find.Execute("54÷8=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=49, 1", 2) | Out-Null
$d.Content.Find.Execute("82÷8=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "90÷4=22, 2", 2) | Out-Null
$d.Content.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2) | Out-Null
$d.Content.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=4, 1", 2) | Out-Null
$d.Content.Find.Execute("18÷4=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "81÷9=9, 0", 2) | Out-Null
$d.Content.Find.Execute("38÷4=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "85÷6=14, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷8=2, 7", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=5, 7", 2) | Out-Null
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "16÷2=8, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷7=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "95÷2=47, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷8=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "62÷7=8, 6", 2) | Out-Null
$d.Content.Find.Execute("39÷8=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=48, 1", 2) | Out-Null
$d.Content.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "45÷2=22, 1", 2) | Out-Null
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "93÷5=18, 3", 2) | Out-Null
